$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# ---- Content Placeholder 2 ("Classification") ----
$shClass = $s.Shapes.Item(2)
$tfClass = $shClass.TextFrame
$tfClass.AutoSize = 2

$trClass = $tfClass.TextRange
$trClass.Text = "Classification:`rThe goal is to predict a categorical or discrete output variable based on input features. `rThe algorithm learns from a labeled dataset where each example is assigned a predefined class or category.`rThe model then uses this training data to classify new, unseen instances into one of the known classes.`rExample: Email spam detection"

# Indent the 4 bullet paragraphs to level 2 (lvl="1")
$trClass.Paragraphs(2).IndentLevel = 2
$trClass.Paragraphs(3).IndentLevel = 2
$trClass.Paragraphs(4).IndentLevel = 2
$trClass.Paragraphs(5).IndentLevel = 2

# ---- Content Placeholder 3 ("Regression") ----
$shReg = $s.Shapes.Item(3)
$tfReg = $shReg.TextFrame
$tfReg.AutoSize = 2

$trReg = $tfReg.TextRange
$trReg.Text = "Regression:`rIt involves predicting a continuous output variable based on input features.`rThe algorithm learns from labeled data, where each example has a corresponding continuous target value.`rThe model then uses this training data to estimate or predict numeric values for new inputs.`rExample: Predicting house prices based on some factors (area, no. of. bedrooms, location...) "

$trReg.Paragraphs(2).IndentLevel = 2
$trReg.Paragraphs(3).IndentLevel = 2
$trReg.Paragraphs(4).IndentLevel = 2
$trReg.Paragraphs(5).IndentLevel = 2
